$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (L) mirroring the existing 2020 column (K).
# L3: empty bottom-border cell matching K3's style
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# L4: year header 2021, formatted like K4 (2020)
$ws.Range("L4").Value = 2021
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)

# L6: Mammals value for 2021
$ws.Range("L6").Value = 7.1
$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)

# L7: Birds value for 2021
$ws.Range("L7").Value = 0.5
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)

# L8: Amphibians/Reptiles value for 2021 ("-" placeholder), styled like K8
$ws.Range("L8").Value = "-"
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)

# Update the active selection to reflect the saved cursor position.
$ws.Range("N5").Select()
